# Update Name of Algo
# Apply numeric corrections to the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 5.885499999999999
$ws.Range("A3").Value  = -21.41010000000002
$ws.Range("B5").Value  = 4.673000000000003
$ws.Range("C5").Value  = -13.9838
$ws.Range("C9").Value  = -11.94320000000001
$ws.Range("C11").Value = -13.27149999999999
$ws.Range("A14").Value = -20.39539999999998
$ws.Range("A21").Value = -21.43450000000001
$ws.Range("C21").Value = -11.27149999999999
$ws.Range("A23").Value = -21.32820000000002
$ws.Range("A25").Value = -22.54890000000004
